$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from the last existing data row (357) down through the new rows (358-366)
# so the new date cells reuse the existing date-style (s="2") instead of minting a new style.
$ws.Range("A357:D357").Copy() | Out-Null
$ws.Range("A358:D366").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the new values (dates 2021-08-24 .. 2021-09-01 and associated counts).
$ws.Cells.Item(358, 1).Value = 44432
$ws.Cells.Item(358, 2).Value = 0
$ws.Cells.Item(358, 3).Value = 4
$ws.Cells.Item(358, 4).Value = 107.0663811563169

$ws.Cells.Item(359, 1).Value = 44433
$ws.Cells.Item(359, 2).Value = 0
$ws.Cells.Item(359, 3).Value = 4
$ws.Cells.Item(359, 4).Value = 107.0663811563169

$ws.Cells.Item(360, 1).Value = 44434
$ws.Cells.Item(360, 2).Value = 1
$ws.Cells.Item(360, 3).Value = 4
$ws.Cells.Item(360, 4).Value = 107.0663811563169

$ws.Cells.Item(361, 1).Value = 44435
$ws.Cells.Item(361, 2).Value = 0
$ws.Cells.Item(361, 3).Value = 3
$ws.Cells.Item(361, 4).Value = 80.29978586723769

$ws.Cells.Item(362, 1).Value = 44436
$ws.Cells.Item(362, 2).Value = 0
$ws.Cells.Item(362, 3).Value = 1
$ws.Cells.Item(362, 4).Value = 26.76659528907923

$ws.Cells.Item(363, 1).Value = 44437
$ws.Cells.Item(363, 2).Value = 0
$ws.Cells.Item(363, 3).Value = 1
$ws.Cells.Item(363, 4).Value = 26.76659528907923

$ws.Cells.Item(364, 1).Value = 44438
$ws.Cells.Item(364, 2).Value = 0
$ws.Cells.Item(364, 3).Value = 1
$ws.Cells.Item(364, 4).Value = 26.76659528907923

$ws.Cells.Item(365, 1).Value = 44439
$ws.Cells.Item(365, 2).Value = 0
$ws.Cells.Item(365, 3).Value = 1
$ws.Cells.Item(365, 4).Value = 26.76659528907923

$ws.Cells.Item(366, 1).Value = 44440
$ws.Cells.Item(366, 2).Value = 0
$ws.Cells.Item(366, 3).Value = 1
$ws.Cells.Item(366, 4).Value = 26.76659528907923
